{"js": "// Replace the three-digit x one-digit multiplication problems in the\n// worksheet table with a new set of problems, preserving formatting.\n// The mapping below is ordered exactly as the expressions occur in the\n// document (row-major through the single table).\nconst replacements = [\n  [\"982\u00d75=\", \"818\u00d74=\"],\n  [\"341\u00d72=\", \"232\u00d76=\"],\n  [\"290\u00d74=\", \"362\u00d73=\"],\n  [\"340\u00d72=\", \"610\u00d77=\"],\n  [\"943\u00d76=\", \"636\u00d78=\"],\n  [\"504\u00d75=\", \"512\u00d76=\"],\n  [\"973\u00d72=\", \"389\u00d74=\"],\n  [\"673\u00d79=\", \"296\u00d76=\"],\n  [\"863\u00d77=\", \"397\u00d74=\"],\n  [\"794\u00d79=\", \"626\u00d73=\"],\n  [\"841\u00d73=\", \"903\u00d78=\"],\n  [\"753\u00d79=\", \"963\u00d73=\"],\n  [\"317\u00d74=\", \"110\u00d78=\"],\n  [\"696\u00d77=\", \"370\u00d74=\"],\n  [\"611\u00d78=\", \"561\u00d75=\"],\n  [\"512\u00d75=\", \"384\u00d74=\"],\n  [\"486\u00d76=\", \"442\u00d79=\"],\n  [\"448\u00d79=\", \"253\u00d76=\"],\n  [\"792\u00d79=\", \"130\u00d79=\"],\n  [\"148\u00d73=\", \"393\u00d79=\"],\n  [\"241\u00d74=\", \"186\u00d76=\"],\n  [\"365\u00d79=\", \"348\u00d79=\"],\n  [\"793\u00d74=\", \"328\u00d79=\"],\n  [\"870\u00d78=\", \"601\u00d72=\"],\n  [\"217\u00d77=\", \"690\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems in the\n# worksheet table with a new set of problems, preserving formatting.\n# The mapping is ordered exactly as the expressions occur in the\n# document (row-major through the single table).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"982\u00d75=\", \"818\u00d74=\"),\n    @(\"341\u00d72=\", \"232\u00d76=\"),\n    @(\"290\u00d74=\", \"362\u00d73=\"),\n    @(\"340\u00d72=\", \"610\u00d77=\"),\n    @(\"943\u00d76=\", \"636\u00d78=\"),\n    @(\"504\u00d75=\", \"512\u00d76=\"),\n    @(\"973\u00d72=\", \"389\u00d74=\"),\n    @(\"673\u00d79=\", \"296\u00d76=\"),\n    @(\"863\u00d77=\", \"397\u00d74=\"),\n    @(\"794\u00d79=\", \"626\u00d73=\"),\n    @(\"841\u00d73=\", \"903\u00d78=\"),\n    @(\"753\u00d79=\", \"963\u00d73=\"),\n    @(\"317\u00d74=\", \"110\u00d78=\"),\n    @(\"696\u00d77=\", \"370\u00d74=\"),\n    @(\"611\u00d78=\", \"561\u00d75=\"),\n    @(\"512\u00d75=\", \"384\u00d74=\"),\n    @(\"486\u00d76=\", \"442\u00d79=\"),\n    @(\"448\u00d79=\", \"253\u00d76=\"),\n    @(\"792\u00d79=\", \"130\u00d79=\"),\n    @(\"148\u00d73=\", \"393\u00d79=\"),\n    @(\"241\u00d74=\", \"186\u00d76=\"),\n    @(\"365\u00d79=\", \"348\u00d79=\"),\n    @(\"793\u00d74=\", \"328\u00d79=\"),\n    @(\"870\u00d78=\", \"601\u00d72=\"),\n    @(\"217\u00d77=\", \"690\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
